$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")
$ws.Activate()

$data = @(
    @(2, 88.5, 11, 46, 90, 7, 0),
    @(3, 65, 21, 53, 133, 7, 0),
    @(4, 50, 11, 42, 71, 7, 0),
    @(5, 60, 14, 48, 53, 7, 0),
    @(6, 85, 23, 55, 119, 7, 0),
    @(7, 75, 10, 57, 154, 7, 0),
    @(8, 37.5, 17, 44, 145, 7, 0),
    @(9, 65, 25, 57, 101, 7, 0),
    @(10, 57, 12, 57, 154, 7, 0),
    @(11, 47.5, 17, 51, 76, 7, 0),
    @(12, 67.5, 10, 41, 54, 7, 0),
    @(13, 85, 13, 50, 81, 7, 0),
    @(14, 57.5, 21, 65, 149, 7, 0),
    @(15, 30, 22, 56, 12, 7, 0),
    @(16, 67.5, 19, 47, 81, 7, 0),
    @(17, 60, 20, 47, 113, 7, 0),
    @(18, 37.5, 12, 56, 172, 7, 0),
    @(19, 57.5, 14, 48, 93, 7, 0),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("L$r").Value = $row[1]
    $ws.Range("M$r").Value = $row[2]
    $ws.Range("N$r").Value = $row[3]
    $ws.Range("O$r").Value = $row[4]
    $ws.Range("P$r").Value = $row[5]
    $ws.Range("Q$r").Value = $row[6]
}

$ws.Range("R2").Formula = "=L2+3*M2+N2+O2"
$ws.Range("R3:R19").Formula = "=L3+3*M3+N3+O3"

$ws.Range("O20").Select()

Write-Output "done"
